$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 93000
$ws.Range("J3").Value = 93000
$ws.Range("L3").Value = 93000
$ws.Range("N3").Value = -93228

$ws.Range("H34").Value = 10753.333
$ws.Range("I34").Value = 7097.625
$ws.Range("K34").Value = 7097.625
$ws.Range("M34").Value = -6894.625

$ws.Range("H36").Value = 10753.333
$ws.Range("I36").Value = 7097.625
$ws.Range("K36").Value = 7097.625
$ws.Range("M36").Value = -6382.625

$ws.Range("H41").Value = 510.33334
$ws.Range("I41").Value = 524.125
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 524.125
$ws.Range("L41").Value = 400
$ws.Range("M41").Value = -84.125
$ws.Range("N41").Value = -1280

$ws.Range("H51").Value = 9421.714
$ws.Range("I51").Value = 8790.6
$ws.Range("K51").Value = 8790.6
$ws.Range("M51").Value = -8306.6

$ws.Range("H74").Value = 6874.75
$ws.Range("I74").Value = 5833.1665
$ws.Range("K74").Value = 5833.1665
$ws.Range("M74").Value = -4897.1665

$ws.Range("H77").Value = 6874.75
$ws.Range("I77").Value = 5833.1665
$ws.Range("K77").Value = 29165.8325
$ws.Range("M77").Value = -24485.8325

$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992

$ws.Range("H102").Value = 93000
$ws.Range("J102").Value = 93000
$ws.Range("L102").Value = 93000
$ws.Range("N102").Value = -99490

$ws.Range("H105").Value = 26459.6
$ws.Range("J105").Value = 26459.6
$ws.Range("L105").Value = 26459.6
$ws.Range("N105").Value = -33447.6

$ws.Range("H111").Value = 998
$ws.Range("I111").Value = 1093.8889
$ws.Range("K111").Value = 3281.6667
$ws.Range("M111").Value = -214.6666999999998

$ws.Range("H132").Value = 3277.0833
$ws.Range("I132").Value = 2674.8223
$ws.Range("J132").Value = 12311
$ws.Range("K132").Value = 8024.466899999999
$ws.Range("L132").Value = 36933
$ws.Range("M132").Value = -5494.466899999999
$ws.Range("N132").Value = -41993

$ws.Range("H137").Value = 9108.857
$ws.Range("I137").Value = 858.6667
$ws.Range("J137").Value = 15296.5
$ws.Range("K137").Value = 2576.0001
$ws.Range("L137").Value = 45889.5
$ws.Range("M137").Value = -26.0001000000002
$ws.Range("N137").Value = -50989.5

$ws.Range("H138").Value = 4771.902
$ws.Range("I138").Value = 2923.5217
$ws.Range("J138").Value = 5890.6577
$ws.Range("K138").Value = 8770.5651
$ws.Range("L138").Value = 17671.9731
$ws.Range("M138").Value = -3630.5651
$ws.Range("N138").Value = -27951.9731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 763.9655
$ws.Range("I2").Value = 622.3333
$ws.Range("J2").Value = 1135.75
$ws.Range("K2").Value = 622.3333
$ws.Range("L2").Value = 1135.75
$ws.Range("M2").Value = -509.3333
$ws.Range("N2").Value = -1361.75

$ws.Range("H32").Value = 26629.607
$ws.Range("I32").Value = 26437.195
$ws.Range("J32").Value = 28399.8
$ws.Range("K32").Value = 26437.195
$ws.Range("L32").Value = 28399.8
$ws.Range("M32").Value = -26150.195
$ws.Range("N32").Value = -28973.8

$ws.Range("H45").Value = 3155.524
$ws.Range("I45").Value = 955.625
$ws.Range("K45").Value = 955.625
$ws.Range("M45").Value = -578.625

$ws.Range("H61").Value = 2126.6365
$ws.Range("J61").Value = 3574.5
$ws.Range("L61").Value = 3574.5
$ws.Range("N61").Value = -3998.5

$ws.Range("H74").Value = 84815
$ws.Range("I74").Value = 112419.664
$ws.Range("J74").Value = 2001
$ws.Range("K74").Value = 112419.664
$ws.Range("L74").Value = 2001
$ws.Range("M74").Value = -111545.664
$ws.Range("N74").Value = -3749

$ws.Range("H77").Value = 84815
$ws.Range("I77").Value = 112419.664
$ws.Range("J77").Value = 2001
$ws.Range("K77").Value = 562098.3200000001
$ws.Range("L77").Value = 10005
$ws.Range("M77").Value = -557730.3200000001
$ws.Range("N77").Value = -18741

$ws.Range("H116").Value = 763.9655
$ws.Range("I116").Value = 622.3333
$ws.Range("J116").Value = 1135.75
$ws.Range("K116").Value = 622.3333
$ws.Range("L116").Value = 1135.75
$ws.Range("M116").Value = 1671.6667
$ws.Range("N116").Value = -5723.75

$ws.Range("H132").Value = 38265.93
$ws.Range("I132").Value = 52492.25
$ws.Range("J132").Value = 6651.8887
$ws.Range("K132").Value = 157476.75
$ws.Range("L132").Value = 19955.6661
$ws.Range("M132").Value = -154946.75
$ws.Range("N132").Value = -25015.6661

$ws.Range("H136").Value = 2126.6365
$ws.Range("J136").Value = 3574.5
$ws.Range("L136").Value = 10723.5
$ws.Range("N136").Value = -15823.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 763.9655
$ws.Range("I3").Value = 622.3333
$ws.Range("J3").Value = 1135.75
$ws.Range("K3").Value = 622.3333
$ws.Range("L3").Value = 1135.75
$ws.Range("M3").Value = -508.3333
$ws.Range("N3").Value = -1363.75

$ws.Range("H75").Value = 9471
$ws.Range("I75").Value = 9471
$ws.Range("K75").Value = 9471
$ws.Range("M75").Value = -8535

$ws.Range("H78").Value = 9471
$ws.Range("I78").Value = 9471
$ws.Range("K78").Value = 28413
$ws.Range("M78").Value = -23733

$ws.Range("H100").Value = 37525.668
$ws.Range("J100").Value = 37525.668
$ws.Range("L100").Value = 37525.668
$ws.Range("N100").Value = -39689.668

$ws.Range("H107").Value = 2166.4285
$ws.Range("I107").Value = 1833
$ws.Range("K107").Value = 1833
$ws.Range("M107").Value = 87

$ws.Range("H134").Value = 1837.5652
$ws.Range("I134").Value = 1746.238
$ws.Range("J134").Value = 2796.5
$ws.Range("K134").Value = 5238.714
$ws.Range("L134").Value = 8389.5
$ws.Range("M134").Value = -2703.714
$ws.Range("N134").Value = -13459.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 805.64703
$ws.Range("I16").Value = 762.25
$ws.Range("K16").Value = 762.25
$ws.Range("M16").Value = -475.25

$ws.Range("H22").Value = 1341.8823
$ws.Range("I22").Value = 309
$ws.Range("K22").Value = 309
$ws.Range("M22").Value = 41

$ws.Range("H31").Value = 4891.3335
$ws.Range("I31").Value = 2789.6365
$ws.Range("J31").Value = 7203.2
$ws.Range("K31").Value = 2789.6365
$ws.Range("L31").Value = 7203.2
$ws.Range("M31").Value = -2494.6365
$ws.Range("N31").Value = -7793.2

$ws.Range("H34").Value = 4891.3335
$ws.Range("I34").Value = 2789.6365
$ws.Range("J34").Value = 7203.2
$ws.Range("K34").Value = 2789.6365
$ws.Range("L34").Value = 7203.2
$ws.Range("M34").Value = -2587.6365
$ws.Range("N34").Value = -7607.2

$ws.Range("H94").Value = 10003174
$ws.Range("I94").Value = 20002770
$ws.Range("K94").Value = 20002770
$ws.Range("M94").Value = -20002319

$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680

$ws.Range("H113").Value = 805.64703
$ws.Range("I113").Value = 762.25
$ws.Range("K113").Value = 762.25
$ws.Range("M113").Value = 1407.75

$ws.Range("H132").Value = 2764.7273
$ws.Range("I132").Value = 2521.2
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 7563.599999999999
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -5033.599999999999
$ws.Range("N132").Value = -20660

$ws.Range("H134").Value = 143940.86
$ws.Range("I134").Value = 167768
$ws.Range("J134").Value = 978
$ws.Range("K134").Value = 503304
$ws.Range("L134").Value = 2934
$ws.Range("M134").Value = -500769
$ws.Range("N134").Value = -8004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 895.7143
$ws.Range("I5").Value = 898.5
$ws.Range("J5").Value = 892
$ws.Range("K5").Value = 2695.5
$ws.Range("L5").Value = 2676
$ws.Range("M5").Value = -2583.5
$ws.Range("N5").Value = -2900

$ws.Range("H21").Value = 501
$ws.Range("I21").Value = 501
$ws.Range("K21").Value = 1503
$ws.Range("M21").Value = -1330

$ws.Range("H68").Value = 452.44446
$ws.Range("I68").Value = 452.44446
$ws.Range("K68").Value = 1357.33338
$ws.Range("M68").Value = -546.33338

$ws.Range("H71").Value = 452.44446
$ws.Range("I71").Value = 452.44446
$ws.Range("K71").Value = 4072.00014
$ws.Range("M71").Value = -16.0001400000001

$ws.Range("H107").Value = 304
$ws.Range("I107").Value = 298
$ws.Range("K107").Value = 894
$ws.Range("M107").Value = 1026

$ws.Range("H108").Value = 4350
$ws.Range("I108").Value = 4350
$ws.Range("K108").Value = 13050
$ws.Range("M108").Value = -10170

$ws.Range("H109").Value = 1112325.4
$ws.Range("I109").Value = 1204
$ws.Range("K109").Value = 3612
$ws.Range("M109").Value = -2572

$ws.Range("H111").Value = 4186.3335
$ws.Range("I111").Value = 780
$ws.Range("J111").Value = 10999
$ws.Range("K111").Value = 2340
$ws.Range("L111").Value = 32997
$ws.Range("M111").Value = 727
$ws.Range("N111").Value = -39131

$ws.Range("H112").Value = 5749.615
$ws.Range("I112").Value = 3474.5
$ws.Range("J112").Value = 13333.333
$ws.Range("K112").Value = 10423.5
$ws.Range("L112").Value = 39999.999
$ws.Range("M112").Value = -9315.5
$ws.Range("N112").Value = -42215.999

$ws.Range("H113").Value = 769.5
$ws.Range("I113").Value = 440.5
$ws.Range("J113").Value = 999.8
$ws.Range("K113").Value = 1321.5
$ws.Range("L113").Value = 2999.4
$ws.Range("M113").Value = 848.5
$ws.Range("N113").Value = -7339.4

$ws.Range("H121").Value = 953.64703
$ws.Range("J121").Value = 1008.2143
$ws.Range("L121").Value = 3024.6429
$ws.Range("N121").Value = -5644.6429

$ws.Range("H123").Value = 2063.375
$ws.Range("I123").Value = 2491.4
$ws.Range("J123").Value = 1350
$ws.Range("K123").Value = 7474.200000000001
$ws.Range("L123").Value = 4050
$ws.Range("M123").Value = -5024.200000000001
$ws.Range("N123").Value = -8950

$ws.Range("H126").Value = 3757.5
$ws.Range("I126").Value = 3757.5
$ws.Range("K126").Value = 11272.5
$ws.Range("M126").Value = -6332.5

$ws.Range("H131").Value = 2331685.8
$ws.Range("J131").Value = 3037508
$ws.Range("L131").Value = 9112524
$ws.Range("N131").Value = -9122604

$ws.Range("H135").Value = 895.7143
$ws.Range("I135").Value = 898.5
$ws.Range("J135").Value = 892
$ws.Range("K135").Value = 8086.5
$ws.Range("L135").Value = 8028
$ws.Range("M135").Value = -5551.5
$ws.Range("N135").Value = -13098

$ws.Range("H137").Value = 4156.391
$ws.Range("I137").Value = 1159.1428
$ws.Range("K137").Value = 3477.4284
$ws.Range("M137").Value = 1622.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.34999999999999
$ws.Range("I2").Value = 84.64286
$ws.Range("K2").Value = 84.64286
$ws.Range("M2").Value = 28.35714

$ws.Range("H20").Value = 1000000
$ws.Range("I20").Value = 1000000
$ws.Range("K20").Value = 1000000
$ws.Range("M20").Value = -999755

$ws.Range("H92").Value = 29900
$ws.Range("I92").Value = 50000
$ws.Range("J92").Value = 9800
$ws.Range("K92").Value = 50000
$ws.Range("L92").Value = 9800
$ws.Range("M92").Value = -48128
$ws.Range("N92").Value = -13544

$ws.Range("H98").Value = 13170.333
$ws.Range("J98").Value = 13170.333
$ws.Range("L98").Value = 13170.333
$ws.Range("N98").Value = -19160.333

$ws.Range("H105").Value = 76249.25
$ws.Range("J105").Value = 76249.25
$ws.Range("L105").Value = 76249.25
$ws.Range("N105").Value = -83237.25

$ws.Range("H107").Value = 67870.60000000001
$ws.Range("I107").Value = 143351.72
$ws.Range("K107").Value = 143351.72
$ws.Range("M107").Value = -141431.72

$ws.Range("H126").Value = 7298.8184
$ws.Range("I126").Value = 6164.5
$ws.Range("K126").Value = 18493.5
$ws.Range("M126").Value = -16023.5

$ws.Range("H132").Value = 103317.4
$ws.Range("I132").Value = 113130.445
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 339391.335
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -336861.335
$ws.Range("N132").Value = -50060

$ws.Range("H134").Value = 43166.668
$ws.Range("J134").Value = 43166.668
$ws.Range("L134").Value = 129500.004
$ws.Range("N134").Value = -134570.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11824.4
$ws.Range("I7").Value = 19742
$ws.Range("K7").Value = 19742
$ws.Range("M7").Value = -19630

$ws.Range("H22").Value = 32296.938
$ws.Range("I22").Value = 45973.047
$ws.Range("J22").Value = 2209.5
$ws.Range("K22").Value = 45973.047
$ws.Range("L22").Value = 2209.5
$ws.Range("M22").Value = -45678.047
$ws.Range("N22").Value = -2799.5

$ws.Range("H27").Value = 32296.938
$ws.Range("I27").Value = 45973.047
$ws.Range("J27").Value = 2209.5
$ws.Range("K27").Value = 45973.047
$ws.Range("L27").Value = 2209.5
$ws.Range("M27").Value = -45866.047
$ws.Range("N27").Value = -2423.5

$ws.Range("H40").Value = 4543.1113
$ws.Range("I40").Value = 4047.5715
$ws.Range("J40").Value = 6277.5
$ws.Range("K40").Value = 4047.5715
$ws.Range("L40").Value = 6277.5
$ws.Range("M40").Value = -3911.5715
$ws.Range("N40").Value = -6549.5

$ws.Range("H46").Value = 12372.454
$ws.Range("I46").Value = 31066.666
$ws.Range("J46").Value = 5362.125
$ws.Range("K46").Value = 31066.666
$ws.Range("L46").Value = 5362.125
$ws.Range("M46").Value = -30878.666
$ws.Range("N46").Value = -5738.125

$ws.Range("H55").Value = 1097
$ws.Range("I55").Value = 981.5263
$ws.Range("J55").Value = 1410.4286
$ws.Range("K55").Value = 981.5263
$ws.Range("L55").Value = 1410.4286
$ws.Range("M55").Value = -808.5263
$ws.Range("N55").Value = -1756.4286

$ws.Range("H68").Value = 3600
$ws.Range("I68").Value = 200
$ws.Range("K68").Value = 200
$ws.Range("M68").Value = 549

$ws.Range("H71").Value = 3600
$ws.Range("I71").Value = 200
$ws.Range("K71").Value = 1000
$ws.Range("M71").Value = 2744

$ws.Range("H82").Value = 2153.8057
$ws.Range("I82").Value = 1269
$ws.Range("J82").Value = 2543.12
$ws.Range("K82").Value = 1269
$ws.Range("L82").Value = 2543.12
$ws.Range("M82").Value = -908
$ws.Range("N82").Value = -3265.12

$ws.Range("H85").Value = 2153.8057
$ws.Range("I85").Value = 1269
$ws.Range("J85").Value = 2543.12
$ws.Range("K85").Value = 1269
$ws.Range("L85").Value = 2543.12
$ws.Range("M85").Value = -21
$ws.Range("N85").Value = -5039.12

$ws.Range("H106").Value = 998
$ws.Range("J106").Value = 998
$ws.Range("L106").Value = 998
$ws.Range("N106").Value = -3522

$ws.Range("H123").Value = 69041
$ws.Range("I123").Value = 34999
$ws.Range("J123").Value = 72445.2
$ws.Range("K123").Value = 34999
$ws.Range("L123").Value = 72445.2
$ws.Range("M123").Value = -30099
$ws.Range("N123").Value = -82245.2

$ws.Range("H126").Value = 11824.4
$ws.Range("I126").Value = 19742
$ws.Range("K126").Value = 59226
$ws.Range("M126").Value = -56756

$ws.Range("H132").Value = 44102.1
$ws.Range("I132").Value = 50394.84
$ws.Range("K132").Value = 151184.52
$ws.Range("M132").Value = -148654.52

$ws.Range("H140").Value = 23762
$ws.Range("J140").Value = 23762
$ws.Range("L140").Value = 23762
$ws.Range("N140").Value = -34122

$ws.Range("H15").Value = 6000
$ws.Range("I15").Value = 6000
$ws.Range("K15").Value = 6000
$ws.Range("M15").Value = -5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 9999
$ws.Range("J18").Value = 9999
$ws.Range("L18").Value = 9999
$ws.Range("N18").Value = -10345

$ws.Range("H81").Value = 2186.2856
$ws.Range("J81").Value = 3444.1667
$ws.Range("L81").Value = 6888.3334
$ws.Range("N81").Value = -9010.3334

$ws.Range("H84").Value = 2186.2856
$ws.Range("J84").Value = 3444.1667
$ws.Range("L84").Value = 34441.667
$ws.Range("N84").Value = -45049.667

$ws.Range("H105").Value = 39600
$ws.Range("J105").Value = 39600
$ws.Range("L105").Value = 39600
$ws.Range("N105").Value = -46588

$ws.Range("H107").Value = 867.1111
$ws.Range("J107").Value = 400
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040

$ws.Range("H112").Value = 29899.5
$ws.Range("J112").Value = 29899.5
$ws.Range("L112").Value = 29899.5
$ws.Range("N112").Value = -32853.5

$ws.Range("H113").Value = 1086.0476
$ws.Range("I113").Value = 478.9
$ws.Range("K113").Value = 1436.7
$ws.Range("M113").Value = 733.3000000000002
